# Update betting odds values on Sheet1 to reflect latest FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 update
$ws.Range("V2").Value = 1.73

# Row 8 updates
$ws.Range("G8").Value  = 2.22
$ws.Range("H8").Value  = 3.35
$ws.Range("I8").Value  = 2.82
$ws.Range("J8").Value  = 2.77
$ws.Range("K8").Value  = 2.2
$ws.Range("L8").Value  = 3.4
$ws.Range("N8").Value  = 8
$ws.Range("O8").Value  = 1.23
$ws.Range("P8").Value  = 3.75
$ws.Range("Q8").Value  = 1.7
$ws.Range("R8").Value  = 2.07
$ws.Range("S8").Value  = 1.34
$ws.Range("T8").Value  = 3
$ws.Range("U8").Value  = 1.57
$ws.Range("V8").Value  = 2.27
$ws.Range("W8").Value  = 10
$ws.Range("Y8").Value  = 8.75
$ws.Range("AA8").Value = 16.5
$ws.Range("AB8").Value = 22
$ws.Range("AC8").Value = 8
$ws.Range("AD8").Value = 6.7
$ws.Range("AG8").Value = 10.75
$ws.Range("AL8").Value = 27
$ws.Range("AO8").Value = 11.5
$ws.Range("AP8").Value = 17
$ws.Range("AR8").Value = 65
$ws.Range("AS8").Value = 175
$ws.Range("AT8").Value = 3
$ws.Range("AU8").Value = 6.5
$ws.Range("AX8").Value = 15.5
$ws.Range("BA8").Value = 90
